$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# New data row (row 47) appended below the existing CRM accuracy log
$ws.Range("A47").Value = 20210811
$ws.Range("B47").Value = 2221.0572992288098
$ws.Range("C47").Value = 2224.4699999999998
$ws.Range("D40:D47").FormulaR1C1 = "=100*(RC[-2]-RC[-1])/RC[-1]"
$ws.Range("E47").Value = 180
$ws.Range("F47").Value = "CRM OPENED 20210720"

# Scroll / selection state, matching where the user left the sheet after entry
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F48").Select()
